$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 161, shifting existing rows 161:212 down to 162:213
$ws.Rows.Item(161).Insert()

# Populate the new row 161 with the new data record
$ws.Cells.Item(161, 1).Value = 7
$ws.Cells.Item(161, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(161, 3).Value = "Ñuble"
$ws.Cells.Item(161, 4).Value = 45027
$ws.Cells.Item(161, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(161, 5).Value = 16
$ws.Cells.Item(161, 6).Value = 100112040
$ws.Cells.Item(161, 7).Value = "Cilantro"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 100
$ws.Cells.Item(161, 11).Value = 1500
$ws.Cells.Item(161, 12).Value = 1500
$ws.Cells.Item(161, 13).Value = 1500
$ws.Cells.Item(161, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(161, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(161, 16).Value = 1500
$ws.Cells.Item(161, 17).Value = 1
$ws.Cells.Item(161, 18).Value = "Hortaliza"
